$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(1, 72861.1282798293),
    @(2, 69335.56608025795),
    @(3, 66537.75603693828),
    @(4, 64167.49999966469),
    @(5, 63812.01336810566),
    @(6, 66555.78069807183),
    @(7, 66451.94826527465),
    @(8, 72289.67388169841),
    @(9, 88683.90869576289),
    @(10, 97359.28029647554),
    @(11, 101156.7777300106),
    @(12, 102750.1808690556),
    @(13, 102937.6130950573),
    @(14, 106470.954737687),
    @(15, 107952.4215135879),
    @(16, 107469.6324529373),
    @(17, 102862.5479607558),
    @(18, 92801.04664529997),
    @(19, 93583.02318242504),
    @(20, 92069.85168083146),
    @(21, 90473.1478324323),
    @(22, 89090.09475814531),
    @(23, 86009.38071636338),
    @(24, 81739.19789414028)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
